$wb = $excel.ActiveWorkbook

$wsValuesets = $wb.Worksheets.Item("Valuesets")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Update the OID values to use the "urn:oid:" prefix (order chosen to match
# the shared-string table insertion order: 4.3.2.1, then 1.2.3.4, then 1.2.3.4.1)
$wsConcepts.Range("D2").Value = "urn:oid:4.3.2.1"
$wsValuesets.Range("B2").Value = "urn:oid:1.2.3.4"
$wsValuesets.Range("B3").Value = "urn:oid:1.2.3.4.1"
$wsConcepts.Range("A2").Value = "urn:oid:1.2.3.4"

# Make the Valuesets sheet the active sheet/tab, with B4 selected
$wsValuesets.Activate()
$wsValuesets.Range("B4").Select()

$wb.Save()
